$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.262.80"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").Value = "3.082.04"
$ws.Range("E3").Value = "  -0.21%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "'562.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").Value = "'144.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.60%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.084.53"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.46%  "

$ws.Range("D10").Value = "'0.155"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.53%  "

$ws.Range("E11").Value = "  -4.43%  "

$ws.Range("D12").Value = "'0.486"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.78%  "

$ws.Range("D13").Value = "'0.0000234"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.82%  "

$ws.Range("D14").Value = "'35.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.31%  "

$ws.Range("D15").Value = "3.574.15"
$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D16").Value = "64.192.58"
$ws.Range("E16").Value = "  -0.59%  "

$ws.Range("D17").Value = "3.075.89"
$ws.Range("E17").Value = "  -0.32%  "

$ws.Range("E18").Value = "  -0.36%  "

$ws.Range("D19").Value = "'6.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.83%  "

$ws.Range("D20").Value = "'491.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.55%  "

$ws.Range("D21").Value = "'14.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.63%  "

$ws.Range("D22").Value = "'14.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.06%  "

$ws.Range("D23").Value = "'0.698"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.03%  "

$ws.Range("D24").Value = "'7.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").Value = "'82.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").Value = "'2.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.17%  "

$ws.Range("D28").Value = "'8.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.63%  "

$ws.Range("D29").Value = "'2.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.13%  "

$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("D31").Value = "'26.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.12%  "

$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("D33").Value = "'2.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.24%  "

$ws.Range("D34").Value = "'5.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.26%  "

$ws.Range("D35").Value = "'6.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.46%  "

$ws.Range("D36").Value = "'55.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.55%  "

$ws.Range("D37").Value = "'0.0415"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.13%  "

$ws.Range("D38").Value = "'446.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.72%  "

$ws.Range("D39").Value = "'0.0825"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.63%  "

$ws.Range("D40").Value = "3.057.23"
$ws.Range("E40").Value = "  +2.68%  "

$ws.Range("D41").Value = "'2.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.03%  "

$ws.Range("D42").Value = "'8.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.19%  "

$ws.Range("D43").Value = "'0.118"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.07%  "

$ws.Range("D44").Value = "'0.282"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.19%  "

$ws.Range("D45").Value = "'28.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").Value = "'2.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.21%  "

$ws.Range("E48").Value = "  +1.34%  "

$ws.Range("D49").Value = "0.0₃0527"
$ws.Range("E49").Value = "  +0.75%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.02%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'117.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.60%  "
